# TOD-E norms run, POM rescale, 24 cell demo strat
#
# 1) Rescale the "ss" lookup column on the existing age-band sheets
#    (5.0-5.3 .. 6.6-6.11, and the first part of 7.0-9.3 which becomes 7.0-7.5).
# 2) Split the old single "7.0-9.3" tab into four narrower age bands:
#    7.0-7.5 (renamed in place), 7.6-7.11, 8.0-8.5, 8.6-9.3 (brand new tabs).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: push a list of B-column ("ss") values into rows 2..N of a sheet,
# leaving column A ("raw") untouched.
# ---------------------------------------------------------------------------
function Set-SSColumn {
    param($Sheet, $Values)
    for ($i = 0; $i -lt $Values.Length; $i++) {
        $Sheet.Cells.Item($i + 2, 2).Value = $Values[$i]
    }
}

# ---------------------------------------------------------------------------
# 1) Update existing sheets 1-6 in place with the new "ss" values.
# ---------------------------------------------------------------------------

$ws1 = $wb.Worksheets.Item(1)   # 5.0-5.3
Set-SSColumn $ws1 @(88,90,92,95,97,100,102,105,107,110,113,115,118,121,124,127,130,130,130,130,130,130,130,130,130)

$ws2 = $wb.Worksheets.Item(2)   # 5.4-5.7
Set-SSColumn $ws2 @(83,85,87,89,92,94,96,99,101,104,107,109,112,115,118,121,124,127,130,130,130,130,130,130,130)

$ws3 = $wb.Worksheets.Item(3)   # 5.8-5.11
Set-SSColumn $ws3 @(78,80,82,85,87,89,91,94,96,99,101,104,107,110,113,116,119,122,125,129,130,130,130,130,130)

$ws4 = $wb.Worksheets.Item(4)   # 6.0-6.5
Set-SSColumn $ws4 @(73,75,77,79,81,84,86,88,91,93,96,98,101,104,107,110,113,117,120,124,127,130,130,130,130)

$ws5 = $wb.Worksheets.Item(5)   # 6.6-6.11
Set-SSColumn $ws5 @(68,70,72,74,76,78,80,82,85,87,90,92,95,98,101,104,108,111,115,119,123,127,130,130,130)

$ws6 = $wb.Worksheets.Item(6)   # was 7.0-9.3, rescaled + renamed 7.0-7.5
Set-SSColumn $ws6 @(64,65,67,69,71,73,75,77,80,82,84,87,90,93,96,99,103,107,112,116,121,126,130,130,130)
$ws6.Name = "7.0-7.5"

# ---------------------------------------------------------------------------
# 2) Append the three brand-new age-band tabs after the renamed sheet,
#    each with a "raw"/"ss" header (bold + centered, matching the other
#    tabs) and 25 data rows.
# ---------------------------------------------------------------------------

function New-NormSheet {
    param($AfterSheet, $Name, $SSValues)

    $ws = $wb.Worksheets.Add($null, $AfterSheet)
    $ws.Name = $Name

    $ws.Range("A1").Value = "raw"
    $ws.Range("B1").Value = "ss"
    $hdr = $ws.Range("A1:B1")
    $hdr.Font.Bold = $true
    $hdr.HorizontalAlignment = -4108

    for ($i = 0; $i -lt $SSValues.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 1).Value = $i + 1
        $ws.Cells.Item($row, 2).Value = $SSValues[$i]
    }

    return $ws
}

$ws7 = New-NormSheet $ws6 "7.6-7.11" @(60,61,63,65,67,69,71,73,75,77,80,82,85,88,92,96,100,105,110,116,121,126,130,130,130)
$ws8 = New-NormSheet $ws7 "8.0-8.5"  @(56,58,59,61,63,65,67,69,71,73,76,78,81,85,88,93,98,104,112,118,123,127,130,130,130)
$ws9 = New-NormSheet $ws8 "8.6-9.3"  @(52,54,55,57,59,60,62,64,66,69,71,74,77,80,85,91,100,111,118,122,125,128,130,130,130)

Write-Host "Final sheet count:" $wb.Worksheets.Count
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    Write-Host $i $wb.Worksheets.Item($i).Name
}
